$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.386.81"
$ws.Range("E2").Value = "  +1.75%  "
$ws.Range("D3").Value = "1.846.15"
$ws.Range("E3").Value = "  +1.48%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.013"
$ws.Range("E4").Value = "  +1.18%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "315.82"
$ws.Range("E5").Value = "  +2.06%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.012"
$ws.Range("E6").Value = "  +1.09%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4742"
$ws.Range("E7").Value = "  +1.41%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3708"
$ws.Range("E8").Value = "  +0.40%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07463"
$ws.Range("E9").Value = "  +1.23%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8879"
$ws.Range("E10").Value = "  +1.86%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.53"
$ws.Range("E11").Value = "  +0.31%  "
$ws.Range("D12").Value = "1.849.03"
$ws.Range("E12").Value = "  +1.42%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07377"
$ws.Range("E13").Value = "  +4.30%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.490"
$ws.Range("E14").Value = "  +2.28%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "93.63"
$ws.Range("E15").Value = "  +1.17%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.598"
$ws.Range("E16").Value = "  +1.48%  "
$ws.Range("E17").Value = "  +1.15%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008850"
$ws.Range("E19").Value = "  +1.05%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "14.88"
$ws.Range("E20").Value = "  +0.69%  "
$ws.Range("D21").Value = "27.400.68"
$ws.Range("E21").Value = "  +1.64%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.342"
$ws.Range("E22").Value = "  -0.14%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "10.74"
$ws.Range("E23").Value = "  +1.48%  "
$ws.Range("D24").Value = "2.071.37"
$ws.Range("E24").Value = "  +1.17%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.908"
$ws.Range("E25").Value = "  +0.30%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "153.06"
$ws.Range("E26").Value = "  +1.07%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.70"
$ws.Range("E27").Value = "  +1.61%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.193"
$ws.Range("E28").Value = "  +0.24%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.309"
$ws.Range("E29").Value = "  -0.36%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "118.22"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08974"
$ws.Range("E31").Value = "  +0.34%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7640"
$ws.Range("E32").Value = "  -0.77%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.181"
$ws.Range("E33").Value = "  +1.28%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.573"
$ws.Range("E34").Value = "  +1.55%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.948"
$ws.Range("E35").Value = "  +1.58%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.013"
$ws.Range("E36").Value = "  +1.17%  "
$ws.Range("E37").Value = "  +2.04%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.05366"
$ws.Range("E38").Value = "  +1.43%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01966"
$ws.Range("E39").Value = "  -0.04%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.005"
$ws.Range("E40").Value = "  +1.89%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.377"
$ws.Range("E41").Value = "  +0.74%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.428"
$ws.Range("E42").Value = "  +2.76%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.5378"
$ws.Range("E43").Value = "  +0.65%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1670"
$ws.Range("E44").Value = "  -0.02%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.589"
$ws.Range("E45").Value = "  +1.61%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4981"
$ws.Range("E46").Value = "  +0.43%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.52"
$ws.Range("E47").Value = "  +0.61%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.013"
$ws.Range("E48").Value = "  +1.23%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.687"
$ws.Range("E49").Value = "  +0.89%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "104.76"
$ws.Range("E50").Value = "  +0.54%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06331"
$ws.Range("E51").Value = "  +0.75%  "
